# Update the "取得日時" (retrieved at) timestamps for the newly appended rows
# from 2026-02-02 12:54:56 to 2026-02-02 13:08:35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2026-02-02 13:08:35"

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
